$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove participant 00106's row (old row 23). Excel's row delete shifts
# everything below up by one and removes the now-unused shared string.
$ws.Rows(23).Delete()

# Participant IDs for the rows that shifted up (now rows 23-31) kept their
# labels but their Age (col B) values were re-entered.
$ws.Cells.Item(23, 2).Value = 66
$ws.Cells.Item(24, 2).Value = 69
$ws.Cells.Item(25, 2).Value = 75
$ws.Cells.Item(26, 2).Value = 72
$ws.Cells.Item(27, 2).Value = 75
$ws.Cells.Item(28, 2).Value = 71
$ws.Cells.Item(29, 2).Value = 61
$ws.Cells.Item(30, 2).Value = 69
$ws.Cells.Item(31, 2).Value = 70
$ws.Cells.Item(2, 3).Value = -0.08195218164326586
$ws.Cells.Item(3, 3).Value = 0.6568106501272286
$ws.Cells.Item(4, 3).Value = 0.6174199209563125
$ws.Cells.Item(5, 3).Value = 0.8433420439841035
$ws.Cells.Item(6, 3).Value = 0.9207427672563853
$ws.Cells.Item(7, 3).Value = -0.8369058374174531
$ws.Cells.Item(8, 3).Value = -0.146261406513147
$ws.Cells.Item(9, 3).Value = -0.2610440497511629
$ws.Cells.Item(10, 3).Value = 0.3219532911962275
$ws.Cells.Item(11, 3).Value = 0.5212777718618362
$ws.Cells.Item(12, 3).Value = 0.7335304552149003
$ws.Cells.Item(13, 3).Value = 0.2892261725922131
$ws.Cells.Item(14, 3).Value = 0.771331611687378
$ws.Cells.Item(15, 3).Value = 0.1334225975328344
$ws.Cells.Item(16, 3).Value = 2.050372167566059
$ws.Cells.Item(17, 3).Value = 0.1795704815443444
$ws.Cells.Item(18, 3).Value = 0.8103745525130024
$ws.Cells.Item(19, 3).Value = 0.6182454285682549
$ws.Cells.Item(20, 3).Value = 0.9386529878247606
$ws.Cells.Item(21, 3).Value = 0.3187844895085998
$ws.Cells.Item(22, 3).Value = 0.7374836884205355
$ws.Cells.Item(23, 3).Value = 0.3802788640284571
$ws.Cells.Item(24, 3).Value = 1.104858743745738
$ws.Cells.Item(25, 3).Value = 1.315290685482777
$ws.Cells.Item(26, 3).Value = 0.6977254268521721
$ws.Cells.Item(27, 3).Value = 1.044920483293541
$ws.Cells.Item(28, 3).Value = 0.6423057088524229
$ws.Cells.Item(29, 3).Value = 0.9471158104941518
$ws.Cells.Item(30, 3).Value = 0.8354543934670168
$ws.Cells.Item(31, 3).Value = 0.4839018293907714
